$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; temporarily unprotect so the cells below can be
# updated, then re-protect it afterwards.
$ws.Unprotect()

# Update the confidential disclaimer date (2021-05-19 -> 2021-05-20)
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-20 for illustrative purposes only and are subject to change."
$ws.Rows(13).AutoFit()

# Update the Weight (column D) and Percent Change (column E) values for rows 2-10
$ws.Range("D2").Value = 0.09077978105910584
$ws.Range("E2").Value = 0.03007715443964964

$ws.Range("D3").Value = 0.105703000141372
$ws.Range("E3").Value = 0.02227020357803844

$ws.Range("D4").Value = 0.1202460034304092
$ws.Range("E4").Value = 0.01149243084726947

$ws.Range("D5").Value = 0.1414161980916658
$ws.Range("E5").Value = 0.01053223990995344

$ws.Range("D6").Value = 0.1379235297958163
$ws.Range("E6").Value = 0.01184999302941558

$ws.Range("D7").Value = 0.1473732465882204
$ws.Range("E7").Value = 0.007910035261602832

$ws.Range("D8").Value = 0.1267904301510683
$ws.Range("E8").Value = 0.01941158629056727

$ws.Range("D9").Value = 0.1297678107423421
$ws.Range("E9").Value = 0.01451475444764028

$ws.Range("D10").Value = 0.9999999999999999
$ws.Range("E10").Value = 0.01510064475904072

# Restore protection on the sheet
$ws.Protect()
